# Update the "vendors" list on the active sheet:
# - Replace the values in A2:A8 with the new agreement/vendor numbers
# - Remove the old row 9 (806538) so the range shrinks from A1:A9 to A1:A8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 672501
$ws.Range("A3").Value = 724629
$ws.Range("A4").Value = 800318
$ws.Range("A5").Value = 801131
$ws.Range("A6").Value = 801254
$ws.Range("A7").Value = 801557
$ws.Range("A8").Value = 806186

# Drop the now-unused last row (previously A9 = 806538)
$ws.Range("A9").Delete()
